# edit.ps1 - Apply the "Saldo" workbook update described by the diff.
#
# Summary of data changes (Export sheet, columns A=Conta, B=Nome, C=Saldo):
#  - Insert 2 rows (ASSAKO 60774.58 / BRUNO 46905.67) before account 005701765
#  - Insert 2 rows (CLISIA 31023.09 / BRUNO 29960.34) before account 004368468
#  - Replace row 004335144/EDMUNDO/10698.6 with two rows:
#      004238436/DIEGO/25076.75 and 004693349/CATARINE/12691.6
#  - Delete row 004890544/ASSAKO/71.86
#  - Insert row 001761119/BLUEMETRIX/63.83 before account 004752519
#  - Delete row 005092207/BRUNO/59.14
#  - Delete row 004551472/DIEGO/43.12
#  - Delete row 004751154/CATARINE/40.6
#  - Delete row 004754056/BRUNO/37.28  (the low-balance duplicate entry)
#  - Insert row 004452597/LARA/36.17 before account 004806286
#  - Delete row 004805269/CLISIA/15.4
#
# NOTE: this runtime's PowerShell function parameter binding only works
# reliably with *positional* arguments (named args like "-row 5" are not
# bound), so all helper functions below are called positionally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Assert-Account($row, $expected) {
    $actual = $ws.Cells.Item($row, 1).Value()
    if ($actual -ne $expected) {
        throw "Row $row mismatch: expected account $expected but found $actual"
    }
}

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = [string]$text
}

function Set-DataRow($row, $conta, $nome, $saldo) {
    Set-TextCell $row 1 $conta
    Set-TextCell $row 2 $nome
    $ws.Cells.Item($row, 3).NumberFormat = "General"
    $ws.Cells.Item($row, 3).Value = [double]$saldo
}

function Insert-DataRowBefore($row, $conta, $nome, $saldo) {
    $ws.Rows.Item($row).Insert()
    Set-DataRow $row $conta $nome $saldo
}

# ---------------------------------------------------------------------
# Apply changes from bottom of the sheet to the top, so row numbers for
# not-yet-processed operations stay valid (rows are 1-based, row 1 is
# the header "Conta/Nome/Saldo"). Each destructive operation is guarded
# by an Assert-Account sanity check against the expected original data.
# ---------------------------------------------------------------------

# Delete row 004805269 / CLISIA / 15.4  (row 136)
Assert-Account 136 "004805269"
$ws.Rows.Item(136).Delete()

# Insert row 004452597 / LARA / 36.17 before account 004806286 (row 99)
Assert-Account 99 "004806286"
Insert-DataRowBefore 99 "004452597" "LARA" 36.17

# Delete row 004754056 / BRUNO / 37.28 (row 97)
Assert-Account 97 "004754056"
$ws.Rows.Item(97).Delete()

# Delete row 004751154 / CATARINE / 40.6 (row 88)
Assert-Account 88 "004751154"
$ws.Rows.Item(88).Delete()

# Delete row 004551472 / DIEGO / 43.12 (row 85)
Assert-Account 85 "004551472"
$ws.Rows.Item(85).Delete()

# Delete row 005092207 / BRUNO / 59.14 (row 68)
Assert-Account 68 "005092207"
$ws.Rows.Item(68).Delete()

# Insert row 001761119 / BLUEMETRIX / 63.83 before account 004752519 (row 63)
Assert-Account 63 "004752519"
Insert-DataRowBefore 63 "001761119" "BLUEMETRIX" 63.83

# Delete row 004890544 / ASSAKO / 71.86 (row 53)
Assert-Account 53 "004890544"
$ws.Rows.Item(53).Delete()

# Replace row 004335144/EDMUNDO/10698.6 (row 6) with two new rows:
#   004238436/DIEGO/25076.75 and 004693349/CATARINE/12691.6
Assert-Account 6 "004335144"
Set-DataRow 6 "004238436" "DIEGO" 25076.75
Insert-DataRowBefore 7 "004693349" "CATARINE" 12691.6

# Insert 2 rows before account 004368468 (row 5):
#   004805273/CLISIA/31023.09 and 004754056/BRUNO/29960.34
Assert-Account 5 "004368468"
Insert-DataRowBefore 5 "004805273" "CLISIA" 31023.09
Insert-DataRowBefore 6 "004754056" "BRUNO" 29960.34

# Insert 2 rows before account 005701765 (row 4):
#   004450724/ASSAKO/60774.58 and 004452912/BRUNO/46905.67
Assert-Account 4 "005701765"
Insert-DataRowBefore 4 "004450724" "ASSAKO" 60774.58
Insert-DataRowBefore 5 "004452912" "BRUNO" 46905.67

Write-Host "Done applying Saldo updates."
